$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "To-Test" (column C) YES/NO flags: move the "YES" marker off
# the rows that have now been tested and onto the next batch of rows.
$ws.Range("C4").Value = "NO"
$ws.Range("C6").Value = "YES"
$ws.Range("C10").Value = "NO"
$ws.Range("C14").Value = "YES"
$ws.Range("C16").Value = "YES"
$ws.Range("C22").Value = "NO"

# Scroll the view down a bit and move the selection to the next cell to
# check (C16), matching where work continues in the sheet.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C16").Select()
